# Update the view-count style numeric counters in column F across the
# four worksheets of the workbook to reflect newly generated stats
# (matches commit "Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

function Set-ColumnFValues($SheetName, $RowToValue) {
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($row in $RowToValue.Keys) {
        $ws.Range("F$row").Value = $RowToValue[$row]
    }
}

# Sheet "展览"
Set-ColumnFValues "展览" @{
    2  = 266
    4  = 9845
    5  = 679
    7  = 353
    8  = 386
    9  = 438
    11 = 217
    13 = 482
    14 = 12478
    18 = 159
    19 = 254
    22 = 129
    23 = 169
    24 = 2745
    30 = 1062
    31 = 4236
    32 = 3743
    33 = 756
    34 = 2643
    36 = 52
    37 = 1348
    38 = 203
    42 = 465
    43 = 609
    45 = 149
    46 = 267
    48 = 142
    49 = 157
}

# Sheet "演出"
Set-ColumnFValues "演出" @{
    8  = 56
    11 = 29
}

# Sheet "本地生活"
Set-ColumnFValues "本地生活" @{
    2 = 59
}

# Sheet "全部类型"
Set-ColumnFValues "全部类型" @{
    4  = 266
    6  = 9845
    7  = 679
    10 = 354
    11 = 386
    12 = 438
    14 = 217
    15 = 482
    16 = 12478
    18 = 59
    19 = 254
    20 = 56
    22 = 129
    23 = 169
    24 = 2745
    29 = 1062
    30 = 4236
    31 = 3743
    32 = 756
    33 = 2643
    35 = 52
    36 = 1348
    37 = 203
    41 = 465
    43 = 609
    45 = 149
    46 = 267
    48 = 142
    49 = 157
}

$wb.Save()
